$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 30.01235008239746
$ws.Range("D2").Value = 0.09235008239745923
$ws.Range("E2").Value = 0.00852853771881751
$ws.Range("C3").Value = 29.95737266540527
$ws.Range("D3").Value = -0.02262733459473054
$ws.Range("E3").Value = 0.0005119962708618896
$ws.Range("B4").Value = 30.03999999999999
$ws.Range("C4").Value = 30.14963340759277
$ws.Range("D4").Value = 0.1096334075927814
$ws.Range("E4").Value = 0.01201948406040494
$ws.Range("B5").Value = 30.21000000000001
$ws.Range("C5").Value = 30.0929012298584
$ws.Range("D5").Value = -0.1170987701416095
$ws.Range("E5").Value = 0.0137121219686775
$ws.Range("C6").Value = 30.26950645446777
$ws.Range("D6").Value = 0.04950645446777457
$ws.Range("E6").Value = 0.002450889033969837
$ws.Range("C7").Value = 30.25530052185059
$ws.Range("D7").Value = -0.1246994781494095
$ws.Range("E7").Value = 0.01554995985073506
$ws.Range("C8").Value = 30.50997734069824
$ws.Range("D8").Value = 0.06997734069824446
$ws.Range("E8").Value = 0.00489682821119818
$ws.Range("C9").Value = 30.38161277770996
$ws.Range("D9").Value = -0.09838722229004304
$ws.Range("E9").Value = 0.009680045509950342
$ws.Range("C10").Value = 30.44911766052246
$ws.Range("D10").Value = -0.2408823394775368
$ws.Range("E10").Value = 0.05802430147217128
$ws.Range("C11").Value = 30.4470043182373
$ws.Range("D11").Value = -0.3029956817626953
$ws.Range("E11").Value = 0.09180638316684053
$ws.Range("C12").Value = 30.63541984558105
$ws.Range("D12").Value = -0.304580154418943
$ws.Range("E12").Value = 0.09276907046586719
$ws.Range("C13").Value = 30.7780818939209
$ws.Range("D13").Value = -0.1719181060791044
$ws.Range("E13").Value = 0.0295558351978262
$ws.Range("C14").Value = 31.21818351745605
$ws.Range("D14").Value = 0.1981835174560587
$ws.Range("E14").Value = 0.03927670659125591
$ws.Range("C15").Value = 31.37577629089355
$ws.Range("D15").Value = 0.2557762908935501
$ws.Range("E15").Value = 0.06542151098326197
$ws.Range("C16").Value = 31.4897289276123
$ws.Range("D16").Value = 0.2097289276123036
$ws.Range("E16").Value = 0.04398622307740686
$ws.Range("C17").Value = 31.24315452575684
$ws.Range("D17").Value = -0.1368454742431595
$ws.Range("E17").Value = 0.01872668382083523
$ws.Range("C18").Value = 31.55834007263184
$ws.Range("D18").Value = -0.02165992736816236
$ws.Range("E18").Value = 0.0004691524535940687
$ws.Range("B19").Value = 31.65000000000001
$ws.Range("C19").Value = 31.96619987487793
$ws.Range("D19").Value = 0.316199874877924
$ws.Range("E19").Value = 0.0999823608728148
$ws.Range("C20").Value = 32.57534408569336
$ws.Range("D20").Value = 0.6953440856933639
$ws.Range("E20").Value = 0.4835033975087402
$ws.Range("C21").Value = 32.38410568237305
$ws.Range("D21").Value = 0.1041056823730457
$ws.Range("E21").Value = 0.01083799310235749
$ws.Range("C22").Value = 32.5141487121582
$ws.Range("D22").Value = 0.06414871215820028
$ws.Range("E22").Value = 0.004115057271555633
$ws.Range("B23").Value = 32.84999999999999
$ws.Range("C23").Value = 32.70607376098633
$ws.Range("D23").Value = -0.1439262390136662
$ws.Range("E23").Value = 0.02071476227661897
$ws.Range("B24").Value = 32.90000000000001
$ws.Range("C24").Value = 32.94353866577148
$ws.Range("D24").Value = 0.04353866577147869
$ws.Range("E24").Value = 0.00189561541716053
$ws.Range("B25").Value = 33.09999999999999
$ws.Range("C25").Value = 32.95751571655273
$ws.Range("D25").Value = -0.1424842834472599
$ws.Range("E25").Value = 0.02030177102947911
$ws.Range("B26").Value = 33.40000000000001
$ws.Range("C26").Value = 33.68034744262695
$ws.Range("D26").Value = 0.2803474426269474
$ws.Range("E26").Value = 0.07859468858746958
$ws.Range("C27").Value = 33.67167282104492
$ws.Range("D27").Value = -0.02832717895508097
$ws.Range("E27").Value = 0.000802429067553182
$ws.Range("B28").Value = 34.09999999999999
$ws.Range("C28").Value = 33.86410903930664
$ws.Range("D28").Value = -0.2358909606933537
$ws.Range("E28").Value = 0.05564454533683334
$ws.Range("B29").Value = 34.40000000000001
$ws.Range("C29").Value = 34.45993423461914
$ws.Range("D29").Value = 0.05993423461913494
$ws.Range("E29").Value = 0.003592112479381513
$ws.Range("B30").Value = 34.90000000000001
$ws.Range("C30").Value = 35.07357406616211
$ws.Range("D30").Value = 0.1735740661621037
$ws.Range("E30").Value = 0.03012795644404635
$ws.Range("C31").Value = 35.75087356567383
$ws.Range("D31").Value = 0.450873565673831
$ws.Range("E31").Value = 0.2032869722234344
$ws.Range("C32").Value = 36.01235580444336
$ws.Range("D32").Value = 0.3123558044433565
$ws.Range("E32").Value = 0.09756614856945639
$ws.Range("C33").Value = 36.00337600708008
$ws.Range("D33").Value = -0.296623992919919
$ws.Range("E33").Value = 0.08798579317575618
$ws.Range("C34").Value = 36.55116653442383
$ws.Range("D34").Value = -0.248833465576169
$ws.Range("E34").Value = 0.0619180935906465
$ws.Range("C35").Value = 37.06188583374023
$ws.Range("D35").Value = -0.2381141662597628
$ws.Range("E35").Value = 0.05669835617358195
$ws.Range("B36").Value = 37.90000000000001
$ws.Range("C36").Value = 37.8508415222168
$ws.Range("D36").Value = -0.04915847778320881
$ws.Range("E36").Value = 0.002416555937962234
$ws.Range("C37").Value = 38.3694953918457
$ws.Range("D37").Value = -0.1305046081542969
$ws.Range("E37").Value = 0.01703145274950657
$ws.Range("B38").Value = 38.90000000000001
$ws.Range("C38").Value = 39.00495529174805
$ws.Range("D38").Value = 0.1049552917480412
$ws.Range("E38").Value = 0.01101561326591644
$ws.Range("B39").Value = 39.40000000000001
$ws.Range("C39").Value = 39.52904891967773
$ws.Range("D39").Value = 0.1290489196777287
$ws.Range("E39").Value = 0.01665362366998887
$ws.Range("B40").Value = 39.90000000000001
$ws.Range("C40").Value = 39.70820999145508
$ws.Range("D40").Value = -0.1917900085449276
$ws.Range("E40").Value = 0.03678340737766338
$ws.Range("B41").Value = 40.09999999999999
$ws.Range("C41").Value = 39.92831802368164
$ws.Range("D41").Value = -0.1716819763183537
$ws.Range("E41").Value = 0.02947470099257576
$ws.Range("B42").Value = 40.59999999999999
$ws.Range("C42").Value = 40.23454666137695
$ws.Range("D42").Value = -0.3654533386230412
$ws.Range("E42").Value = 0.1335561427107272
$ws.Range("B43").Value = 40.90000000000001
$ws.Range("C43").Value = 40.47024154663086
$ws.Range("D43").Value = -0.4297584533691463
$ws.Range("E43").Value = 0.1846923282422407
$ws.Range("B44").Value = 41.20000000000001
$ws.Range("C44").Value = 41.21378707885742
$ws.Range("D44").Value = 0.01378707885741193
$ws.Range("E44").Value = 0.000190083543420495
$ws.Range("C45").Value = 41.2248649597168
$ws.Range("D45").Value = -0.2751350402832031
$ws.Range("E45").Value = 0.07569929039163981
$ws.Range("C46").Value = 41.653076171875
$ws.Range("D46").Value = -0.1469238281249972
$ws.Range("E46").Value = 0.0215866112709037
$ws.Range("C47").Value = 42.15799331665039
$ws.Range("D47").Value = -0.04200668334961222
$ws.Range("E47").Value = 0.001764561446034588
$ws.Range("C48").Value = 43.88340759277344
$ws.Range("D48").Value = 1.183407592773435
$ws.Range("E48").Value = 1.400453530633815
$ws.Range("C49").Value = 44.43407821655273
$ws.Range("D49").Value = 0.7340782165527386
$ws.Range("E49").Value = 0.5388708280172494
$ws.Range("C50").Value = 44.31099700927734
$ws.Range("D50").Value = 0.1109970092773409
$ws.Range("E50").Value = 0.0123203360685141
$ws.Range("C51").Value = 44.46113967895508
$ws.Range("D51").Value = -1.138860321044923
$ws.Range("E51").Value = 1.297002830850546
$ws.Range("C52").Value = -0.05531524658206166
$ws.Range("E52").Value = 5.604465680149265
$ws.Range("E53").Value = 0.1120893136029853
